$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 50000
$ws.Range("I13").Value = 50000
$ws.Range("K13").Value = 50000
$ws.Range("M13").Value = -49831

$ws.Range("H19").Value = 710.3103599999999
$ws.Range("I19").Value = 595.5625
$ws.Range("J19").Value = 851.53845
$ws.Range("K19").Value = 595.5625
$ws.Range("L19").Value = 851.53845
$ws.Range("M19").Value = -420.5625
$ws.Range("N19").Value = -1201.53845

$ws.Range("H33").Value = 49
$ws.Range("I33").Value = 49
$ws.Range("K33").Value = 49
$ws.Range("M33").Value = 180

$ws.Range("H76").Value = 4276496
$ws.Range("I76").Value = 4833787.5
$ws.Range("K76").Value = 4833787.5
$ws.Range("M76").Value = -4833472.5

$ws.Range("H79").Value = 4276496
$ws.Range("I79").Value = 4833787.5
$ws.Range("K79").Value = 4833787.5
$ws.Range("M79").Value = -4832695.5

$ws.Range("H132").Value = 282741.9
$ws.Range("I132").Value = 369148.28
$ws.Range("K132").Value = 1107444.84
$ws.Range("M132").Value = -1104914.84

$ws.Range("H135").Value = 1198.3784
$ws.Range("I135").Value = 1037.5758
$ws.Range("J135").Value = 2525
$ws.Range("K135").Value = 9338.182200000001
$ws.Range("L135").Value = 22725
$ws.Range("M135").Value = -6803.182200000001
$ws.Range("N135").Value = -27795

$ws.Range("H136").Value = 57166.668
$ws.Range("J136").Value = 57166.668
$ws.Range("L136").Value = 57166.668
$ws.Range("N136").Value = -67366.66800000001

$ws.Range("H138").Value = 5448870
$ws.Range("I138").Value = 1852493.2
$ws.Range("J138").Value = 6947360.5
$ws.Range("K138").Value = 5557479.6
$ws.Range("L138").Value = 20842081.5
$ws.Range("M138").Value = -5552339.6
$ws.Range("N138").Value = -20852361.5


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19426.666
$ws.Range("I32").Value = 2572.9644
$ws.Range("K32").Value = 2572.9644
$ws.Range("M32").Value = -2285.9644

$ws.Range("H45").Value = 882.4
$ws.Range("I45").Value = 737.3333
$ws.Range("K45").Value = 737.3333
$ws.Range("M45").Value = -360.3333

$ws.Range("H61").Value = 3385.4443
$ws.Range("I61").Value = 2704.1667
$ws.Range("K61").Value = 2704.1667
$ws.Range("M61").Value = -2492.1667

$ws.Range("H136").Value = 3385.4443
$ws.Range("I136").Value = 2704.1667
$ws.Range("K136").Value = 8112.500100000001
$ws.Range("M136").Value = -5562.500100000001


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 40000
$ws.Range("J55").Value = 40000
$ws.Range("L55").Value = 40000
$ws.Range("N55").Value = -40546

$ws.Range("H94").Value = 1316.091
$ws.Range("I94").Value = 1511.875
$ws.Range("J94").Value = 794
$ws.Range("K94").Value = 1511.875
$ws.Range("L94").Value = 794
$ws.Range("M94").Value = -1060.875
$ws.Range("N94").Value = -1696

$ws.Range("H99").Value = 2034.6154
$ws.Range("I99").Value = 1890
$ws.Range("J99").Value = 2125
$ws.Range("K99").Value = 1890
$ws.Range("L99").Value = 2125
$ws.Range("M99").Value = -392
$ws.Range("N99").Value = -5121

$ws.Range("H107").Value = 936.5
$ws.Range("I107").Value = 698
$ws.Range("J107").Value = 1493
$ws.Range("K107").Value = 698
$ws.Range("L107").Value = 1493
$ws.Range("M107").Value = 1222
$ws.Range("N107").Value = -5333


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1040.3103
$ws.Range("I31").Value = 934.6070999999999
$ws.Range("K31").Value = 934.6070999999999
$ws.Range("M31").Value = -639.6070999999999

$ws.Range("H34").Value = 1040.3103
$ws.Range("I34").Value = 934.6070999999999
$ws.Range("K34").Value = 934.6070999999999
$ws.Range("M34").Value = -732.6070999999999

$ws.Range("H58").Value = 3230.5625
$ws.Range("I58").Value = 3136
$ws.Range("J58").Value = 3273.5454
$ws.Range("K58").Value = 3136
$ws.Range("L58").Value = 3273.5454
$ws.Range("M58").Value = -2933
$ws.Range("N58").Value = -3679.5454

$ws.Range("H132").Value = 2739.8215
$ws.Range("I132").Value = 2254.0952
$ws.Range("J132").Value = 4197
$ws.Range("K132").Value = 6762.285600000001
$ws.Range("L132").Value = 12591
$ws.Range("M132").Value = -4232.285600000001
$ws.Range("N132").Value = -17651

$ws.Range("H136").Value = 3230.5625
$ws.Range("I136").Value = 3136
$ws.Range("J136").Value = 3273.5454
$ws.Range("K136").Value = 9408
$ws.Range("L136").Value = 9820.636200000001
$ws.Range("M136").Value = -6858
$ws.Range("N136").Value = -14920.6362


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 11364401
$ws.Range("J113").Value = 13889691
$ws.Range("L113").Value = 41669073
$ws.Range("N113").Value = -41673413

$ws.Range("H132").Value = 1241.4762
$ws.Range("I132").Value = 826.3333
$ws.Range("J132").Value = 1407.5333
$ws.Range("K132").Value = 7436.9997
$ws.Range("L132").Value = 12667.7997
$ws.Range("M132").Value = -4906.9997
$ws.Range("N132").Value = -17727.7997

$ws.Range("H136").Value = 2663.0605
$ws.Range("I136").Value = 1700.25
$ws.Range("J136").Value = 2795.862
$ws.Range("K136").Value = 5100.75
$ws.Range("L136").Value = 8387.585999999999
$ws.Range("M136").Value = -0.75
$ws.Range("N136").Value = -18587.586

$ws.Range("H139").Value = 2470.45
$ws.Range("I139").Value = 2024.0588
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 6072.1764
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -932.1764000000003
$ws.Range("N139").Value = -25280


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5285.2915
$ws.Range("I70").Value = 5340.5137
$ws.Range("J70").Value = 5099.5454
$ws.Range("K70").Value = 5340.5137
$ws.Range("L70").Value = 5099.5454
$ws.Range("M70").Value = -5070.5137
$ws.Range("N70").Value = -5639.5454

$ws.Range("H73").Value = 5285.2915
$ws.Range("I73").Value = 5340.5137
$ws.Range("J73").Value = 5099.5454
$ws.Range("K73").Value = 5340.5137
$ws.Range("L73").Value = 5099.5454
$ws.Range("M73").Value = -4404.5137
$ws.Range("N73").Value = -6971.5454

$ws.Range("H122").Value = 1390314.8
$ws.Range("I122").Value = 1853003
$ws.Range("K122").Value = 5559009
$ws.Range("M122").Value = -5556559

$ws.Range("H132").Value = 5881.25
$ws.Range("I132").Value = 6702
$ws.Range("J132").Value = 4513.3335
$ws.Range("K132").Value = 20106
$ws.Range("L132").Value = 13540.0005
$ws.Range("M132").Value = -17576
$ws.Range("N132").Value = -18600.0005


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2895.6086
$ws.Range("I7").Value = 1499.8572
$ws.Range("K7").Value = 1499.8572
$ws.Range("M7").Value = -1387.8572

$ws.Range("H16").Value = 22698
$ws.Range("I16").Value = 3372.75
$ws.Range("J16").Value = 99999
$ws.Range("K16").Value = 3372.75
$ws.Range("L16").Value = 99999
$ws.Range("M16").Value = -3202.75
$ws.Range("N16").Value = -100339

$ws.Range("H100").Value = 2091.2285
$ws.Range("I100").Value = 1292.6
$ws.Range("J100").Value = 2690.2
$ws.Range("K100").Value = 1292.6
$ws.Range("L100").Value = 2690.2
$ws.Range("M100").Value = -751.5999999999999
$ws.Range("N100").Value = -3772.2

$ws.Range("H122").Value = 3068.6
$ws.Range("I122").Value = 1936.3636
$ws.Range("J122").Value = 3587.5417
$ws.Range("K122").Value = 5809.0908
$ws.Range("L122").Value = 10762.6251
$ws.Range("M122").Value = -3359.0908
$ws.Range("N122").Value = -15662.6251

$ws.Range("H126").Value = 2895.6086
$ws.Range("I126").Value = 1499.8572
$ws.Range("K126").Value = 4499.571599999999
$ws.Range("M126").Value = -2029.571599999999

$ws.Range("H132").Value = 4565.026
$ws.Range("I132").Value = 4343.6
$ws.Range("J132").Value = 4703.4165
$ws.Range("K132").Value = 13030.8
$ws.Range("L132").Value = 14110.2495
$ws.Range("M132").Value = -10500.8
$ws.Range("N132").Value = -19170.2495

$ws.Range("H136").Value = 4382.273
$ws.Range("I136").Value = 2043.826
$ws.Range("J136").Value = 9760.700000000001
$ws.Range("K136").Value = 6131.478
$ws.Range("L136").Value = 29282.1
$ws.Range("M136").Value = -3581.478
$ws.Range("N136").Value = -34382.10000000001


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 907.2308
$ws.Range("I113").Value = 757.3
$ws.Range("K113").Value = 2271.9
$ws.Range("M113").Value = -101.8999999999996

$ws.Range("H132").Value = 14289484
$ws.Range("I132").Value = 21743352
$ws.Range("J132").Value = 2905
$ws.Range("K132").Value = 65230056
$ws.Range("L132").Value = 8715
$ws.Range("M132").Value = -65227526
$ws.Range("N132").Value = -13775

$ws.Range("H136").Value = 11942598
$ws.Range("I136").Value = 20897360
$ws.Range("J136").Value = 2915.1667
$ws.Range("K136").Value = 62692080
$ws.Range("L136").Value = 8745.500100000001
$ws.Range("M136").Value = -62689530
$ws.Range("N136").Value = -13845.5001

